$d = $word.ActiveDocument

# --- 1. "unknown unknown" -> "unknown unknowns" + new sentence about
#        "In the context of economic policy, unknown unknowns" replacing
#        "These are" -------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "unknown unknown (Rumsfeld, 2011).  These are the gaps",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "unknown unknowns (Rumsfeld, 2011).  In the context of economic policy, unknown unknowns are the gaps",
    2)

# --- 2. "mental constructs" -> "conceptualizations" ------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "previous theories and mental constructs less accurate",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "previous theories and conceptualizations less accurate",
    2)

# --- 3. Remove the _GoBack bookmark from the end of the Response
#        paragraph; it gets re-created further down after the new
#        Reference paragraph is built. ---------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 4. Split the References paragraph: "New York, NY: Penguin Group."
#        becomes its own paragraph end; the trailing two blank-space runs
#        move into a brand-new trailing paragraph. ------------------------
$r3 = $d.Content
$r3.Find.Execute(
    "New York, NY: Penguin Group.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.InsertParagraphAfter()

# --- 5. The new trailing paragraph now holds two space characters
#        (one from each of the old trailing space runs); drop the first
#        one so only a single space run remains. ---------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$dropSpace = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + 1)
$dropSpace.Delete()

# --- 6. Re-create the _GoBack bookmark at the very start of that
#        trailing paragraph. -------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bmRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 7. Italicize the book title within the reference entry. ---------------
$r4 = $d.Content
$r4.Find.Execute(
    "Known and unknown: A memoir",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.Font.Italic = $true
